# Update "想去人数" (want-to-go count) figures in column F for the
# "展览" (Exhibitions) sheet and the "全部类型" (All types) sheet,
# reflecting newly generated stats output.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (sheet1) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3052
$ws1.Range("F3").Value = 473
$ws1.Range("F4").Value = 53
$ws1.Range("F7").Value = 1035
$ws1.Range("F8").Value = 14684
$ws1.Range("F9").Value = 172
$ws1.Range("F10").Value = 132
$ws1.Range("F11").Value = 5855
$ws1.Range("F13").Value = 81
$ws1.Range("F14").Value = 47
$ws1.Range("F15").Value = 73
$ws1.Range("F18").Value = 89
$ws1.Range("F19").Value = 190
$ws1.Range("F20").Value = 805
$ws1.Range("F21").Value = 2942
$ws1.Range("F22").Value = 84
$ws1.Range("F23").Value = 10650
$ws1.Range("F24").Value = 1205
$ws1.Range("F25").Value = 69
$ws1.Range("F26").Value = 98
$ws1.Range("F27").Value = 3743
$ws1.Range("F28").Value = 248

# ---- Sheet "全部类型" (sheet4) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3052
$ws4.Range("F4").Value = 473
$ws4.Range("F5").Value = 53
$ws4.Range("F8").Value = 1035
$ws4.Range("F9").Value = 14684
$ws4.Range("F10").Value = 172
$ws4.Range("F11").Value = 132
$ws4.Range("F12").Value = 5855
$ws4.Range("F14").Value = 81
$ws4.Range("F15").Value = 47
$ws4.Range("F16").Value = 73
$ws4.Range("F19").Value = 89
$ws4.Range("F20").Value = 190
$ws4.Range("F21").Value = 805
$ws4.Range("F22").Value = 2942
$ws4.Range("F23").Value = 84
$ws4.Range("F25").Value = 10650
$ws4.Range("F26").Value = 1205
$ws4.Range("F27").Value = 69
$ws4.Range("F28").Value = 98
$ws4.Range("F29").Value = 3743
$ws4.Range("F30").Value = 248
